# Apply the commit "Add files via upload":
#   1. Remove the "Team Members" slide (the 2nd slide in the deck).
#   2. Refresh the cached text of every "datetimeFigureOut" date
#      placeholder (slide master, the 5 slide layouts, the notes master
#      and the handout master) from 12/6/2020 to 8/9/2021.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Delete the "Team Members" slide.
# ---------------------------------------------------------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isTeamMembers = $false
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "Team Members") {
                $isTeamMembers = $true
            }
        }
    }
    if ($isTeamMembers) {
        $slide.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Update the cached "datetimeFigureOut" placeholder text everywhere
#    it appears: slide master, every slide layout, notes master, and
#    handout master.
# ---------------------------------------------------------------------
$oldDate = "12/6/2020"
$newDate = "8/9/2021"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom layout hanging off the (single) design/master.
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $design = $p.Designs.Item($d)
    $layouts = $design.SlideMaster.CustomLayouts
    for ($L = 1; $L -le $layouts.Count; $L++) {
        Update-DatePlaceholder $layouts.Item($L).Shapes
    }
}

# Notes master
if ($p.HasNotesMaster) {
    Update-DatePlaceholder $p.NotesMaster.Shapes
} else {
    Update-DatePlaceholder $p.NotesMaster.Shapes
}

# Handout master
if ($p.HasHandoutMaster) {
    Update-DatePlaceholder $p.HandoutMaster.Shapes
} else {
    Update-DatePlaceholder $p.HandoutMaster.Shapes
}
